$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (reflects the new "through" date)
$ws.Name = "Through 2022-02-18"

# Update the label for the February row
$ws.Range("A3").Value = "February (through 02-18)"

# Update February row values (row 3)
$ws.Range("E3").Value = 35
$ws.Range("F3").Value = 19
$ws.Range("G3").Value = 43
$ws.Range("H3").Value = 80
$ws.Range("I3").Value = 85

# Update Total row values (row 4)
$ws.Range("E4").Value = 121
$ws.Range("F4").Value = 68
$ws.Range("G4").Value = 117
$ws.Range("H4").Value = 297
$ws.Range("I4").Value = 245
